$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Atualiza os dados de faturamento anual para o ano de 2025 (linha 7)
$ws.Range("B7").Value = 2937387.51
$ws.Range("C7").Value = -33.88845681605977
$ws.Range("D7").Value = 2970
$ws.Range("E7").Value = 2970
$ws.Range("F7").Value = 989.0193636363636
$ws.Range("G7").Value = 5.422312632707382
